$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 971.25
$ws.Range("J7").Value = 971.25
$ws.Range("L7").Value = 971.25
$ws.Range("N7").Value = -1195.25
$ws.Range("H14").Value = 971.25
$ws.Range("J14").Value = 971.25
$ws.Range("L14").Value = 971.25
$ws.Range("N14").Value = -1353.25
$ws.Range("H17").Value = 1180.8334
$ws.Range("J17").Value = 1180.8334
$ws.Range("L17").Value = 3542.5002
$ws.Range("N17").Value = -3878.5002
$ws.Range("H19").Value = 563.8421
$ws.Range("I19").Value = 323.625
$ws.Range("J19").Value = 738.5454999999999
$ws.Range("K19").Value = 323.625
$ws.Range("L19").Value = 738.5454999999999
$ws.Range("M19").Value = -148.625
$ws.Range("N19").Value = -1088.5455
$ws.Range("H64").Value = 4946
$ws.Range("I64").Value = 5200
$ws.Range("J64").Value = 4800.857
$ws.Range("K64").Value = 5200
$ws.Range("L64").Value = 4800.857
$ws.Range("M64").Value = -4952
$ws.Range("N64").Value = -5296.857
$ws.Range("H67").Value = 4946
$ws.Range("I67").Value = 5200
$ws.Range("J67").Value = 4800.857
$ws.Range("K67").Value = 5200
$ws.Range("L67").Value = 4800.857
$ws.Range("M67").Value = -4342
$ws.Range("N67").Value = -6516.857
$ws.Range("H93").Value = 66000
$ws.Range("J93").Value = 66000
$ws.Range("L93").Value = 66000
$ws.Range("N93").Value = -70992
$ws.Range("H106").Value = 3448.077
$ws.Range("I106").Value = 2853.125
$ws.Range("J106").Value = 4400
$ws.Range("K106").Value = 2853.125
$ws.Range("L106").Value = 4400
$ws.Range("M106").Value = -2222.125
$ws.Range("N106").Value = -5662
$ws.Range("H113").Value = 3361.2273
$ws.Range("J113").Value = 3898.4
$ws.Range("L113").Value = 3898.4
$ws.Range("N113").Value = -10406.4
$ws.Range("H115").Value = 3420.818
$ws.Range("I115").Value = 661.5714
$ws.Range("J115").Value = 8249.5
$ws.Range("K115").Value = 1984.7142
$ws.Range("L115").Value = 24748.5
$ws.Range("M115").Value = -417.7142000000001
$ws.Range("N115").Value = -27882.5
$ws.Range("H118").Value = 1416.9231
$ws.Range("I118").Value = 713.3333
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 2139.9999
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = -482.9998999999998
$ws.Range("N118").Value = -12314
$ws.Range("H137").Value = 1788161.5
$ws.Range("I137").Value = 3031952.8
$ws.Range("K137").Value = 9095858.399999999
$ws.Range("M137").Value = -9093308.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 120.166664
$ws.Range("I5").Value = 120.166664
$ws.Range("K5").Value = 120.166664
$ws.Range("M5").Value = -8.166663999999997
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 500
$ws.Range("M16").Value = -213
$ws.Range("H32").Value = 11003.41
$ws.Range("I32").Value = 7901
$ws.Range("J32").Value = 19833.346
$ws.Range("K32").Value = 7901
$ws.Range("L32").Value = 19833.346
$ws.Range("M32").Value = -7614
$ws.Range("N32").Value = -20407.346

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 120.166664
$ws.Range("I4").Value = 120.166664
$ws.Range("K4").Value = 120.166664
$ws.Range("M4").Value = -5.166663999999997
$ws.Range("H103").Value = 49900
$ws.Range("J103").Value = 49900
$ws.Range("L103").Value = 49900
$ws.Range("N103").Value = -52244

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 303048.2
$ws.Range("J31").Value = 462002.9
$ws.Range("L31").Value = 462002.9
$ws.Range("N31").Value = -462592.9
$ws.Range("H34").Value = 303048.2
$ws.Range("J34").Value = 462002.9
$ws.Range("L34").Value = 462002.9
$ws.Range("N34").Value = -462406.9
$ws.Range("H58").Value = 26317508
$ws.Range("I58").Value = 37038530
$ws.Range("J58").Value = 2275.182
$ws.Range("K58").Value = 37038530
$ws.Range("L58").Value = 2275.182
$ws.Range("M58").Value = -37038327
$ws.Range("N58").Value = -2681.182
$ws.Range("H132").Value = 47656.91
$ws.Range("I132").Value = 1853.9286
$ws.Range("J132").Value = 127812.125
$ws.Range("K132").Value = 5561.7858
$ws.Range("L132").Value = 383436.375
$ws.Range("M132").Value = -3031.7858
$ws.Range("N132").Value = -388496.375
$ws.Range("H134").Value = 61091.5
$ws.Range("I134").Value = 870.6667
$ws.Range("J134").Value = 73135.664
$ws.Range("K134").Value = 2612.0001
$ws.Range("L134").Value = 219406.992
$ws.Range("M134").Value = -77.0001000000002
$ws.Range("N134").Value = -224476.992
$ws.Range("H136").Value = 26317508
$ws.Range("I136").Value = 37038530
$ws.Range("J136").Value = 2275.182
$ws.Range("K136").Value = 111115590
$ws.Range("L136").Value = 6825.545999999999
$ws.Range("M136").Value = -111113040
$ws.Range("N136").Value = -11925.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 857.1429000000001
$ws.Range("I17").Value = 875
$ws.Range("J17").Value = 833.3333
$ws.Range("K17").Value = 2625
$ws.Range("L17").Value = 2499.9999
$ws.Range("M17").Value = -2456
$ws.Range("N17").Value = -2837.9999
$ws.Range("H104").Value = 3997.5
$ws.Range("J104").Value = 3997.5
$ws.Range("L104").Value = 11992.5
$ws.Range("N104").Value = -17234.5
$ws.Range("H131").Value = 931.6842
$ws.Range("J131").Value = 1105.2142
$ws.Range("L131").Value = 3315.6426
$ws.Range("N131").Value = -13395.6426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 126.71429
$ws.Range("I2").Value = 90.85714
$ws.Range("J2").Value = 162.57143
$ws.Range("K2").Value = 90.85714
$ws.Range("L2").Value = 162.57143
$ws.Range("M2").Value = 22.14286
$ws.Range("N2").Value = -388.57143
$ws.Range("H132").Value = 103443.75
$ws.Range("I132").Value = 93305.73
$ws.Range("J132").Value = 115834.664
$ws.Range("K132").Value = 279917.19
$ws.Range("L132").Value = 347503.992
$ws.Range("M132").Value = -277387.19
$ws.Range("N132").Value = -352563.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 977.4783
$ws.Range("I16").Value = 769.5
$ws.Range("J16").Value = 2364
$ws.Range("K16").Value = 769.5
$ws.Range("L16").Value = 2364
$ws.Range("M16").Value = -599.5
$ws.Range("N16").Value = -2704
$ws.Range("H22").Value = 685.76
$ws.Range("I22").Value = 599.1818
$ws.Range("J22").Value = 753.7857
$ws.Range("K22").Value = 599.1818
$ws.Range("L22").Value = 753.7857
$ws.Range("M22").Value = -304.1818
$ws.Range("N22").Value = -1343.7857
$ws.Range("H27").Value = 685.76
$ws.Range("I27").Value = 599.1818
$ws.Range("J27").Value = 753.7857
$ws.Range("K27").Value = 599.1818
$ws.Range("L27").Value = 753.7857
$ws.Range("M27").Value = -492.1818
$ws.Range("N27").Value = -967.7857
$ws.Range("H82").Value = 2669
$ws.Range("I82").Value = 2001
$ws.Range("K82").Value = 2001
$ws.Range("M82").Value = -1640
$ws.Range("H85").Value = 2669
$ws.Range("I85").Value = 2001
$ws.Range("K85").Value = 2001
$ws.Range("M85").Value = -753
$ws.Range("H93").Value = 1168.4814
$ws.Range("I93").Value = 1154.7142
$ws.Range("K93").Value = 1154.7142
$ws.Range("M93").Value = 93.28580000000011
$ws.Range("H100").Value = 1822
$ws.Range("I100").Value = 1516
$ws.Range("J100").Value = 1975
$ws.Range("K100").Value = 1516
$ws.Range("L100").Value = 1975
$ws.Range("M100").Value = -975
$ws.Range("N100").Value = -3057
$ws.Range("H122").Value = 4093.9565
$ws.Range("I122").Value = 5885.857
$ws.Range("K122").Value = 17657.571
$ws.Range("M122").Value = -15207.571
$ws.Range("H132").Value = 69937
$ws.Range("I132").Value = 6165.5835
$ws.Range("J132").Value = 261251.25
$ws.Range("K132").Value = 18496.7505
$ws.Range("L132").Value = 783753.75
$ws.Range("M132").Value = -15966.7505
$ws.Range("N132").Value = -788813.75
$ws.Range("H136").Value = 73315.89999999999
$ws.Range("I136").Value = 42736.96
$ws.Range("J136").Value = 338333.34
$ws.Range("K136").Value = 128210.88
$ws.Range("L136").Value = 1015000.02
$ws.Range("M136").Value = -125660.88
$ws.Range("N136").Value = -1020100.02

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 695.1613
$ws.Range("I113").Value = 857.05554
$ws.Range("J113").Value = 471
$ws.Range("K113").Value = 2571.16662
$ws.Range("L113").Value = 1413
$ws.Range("M113").Value = -401.16662
$ws.Range("N113").Value = -5753
$ws.Range("H136").Value = 54076.36
$ws.Range("I136").Value = 54381.74
$ws.Range("J136").Value = 53786.25
$ws.Range("K136").Value = 163145.22
$ws.Range("L136").Value = 161358.75
$ws.Range("M136").Value = -160595.22
$ws.Range("N136").Value = -166458.75
